# "something is wrong with irrigation electification"
#
# A missing electricity conversion-factor row is inserted into the
# "energy" sheet (between the EPA-sourced mj/l/mj/kwh rows and the
# Grassini & Cassman rows), sourced from Wikipedia's
# "Miles per gallon gasoline equivalent" article. The previously blank
# formatting left on column E (an inert "apply fill" style with no
# actual fill) is cleaned up while we're in there, and a stray blank
# styled cell on the "manufacture" sheet is cleared too.

$wb = $excel.ActiveWorkbook

# --- "energy" sheet: insert the missing electricity row -------------------
$ws1 = $wb.Worksheets.Item("energy")

# Row 9 was: diesel | Grassini and Cassman 2012 | mj/l | 43 | Table S3
# Insert a fresh row above it for the missing electricity entry, pushing
# the diesel/gasoline/electricity (Grassini) and diesel/gasoline/electricity
# (Hoffman) rows down by one (old rows 9-14 -> new rows 10-15).
$ws1.Rows.Item(9).Insert()

$ws1.Range("A9").Value = "electricity"
$ws1.Range("B9").Value = "EPA"
$ws1.Range("C9").Value = "mj/kwh"
$ws1.Range("D9").Value = 3.6
$ws1.Range("E9").Value = "https://en.wikipedia.org/wiki/Miles_per_gallon_gasoline_equivalent"

# Clean up the inert leftover style (fontId/fill/border all default -- it
# renders identically to no style at all) that was sitting on column E.
$ws1.Range("E7:E15").ClearFormats()

# --- "manufacture" sheet: drop the stray empty styled cell ----------------
$ws3 = $wb.Worksheets.Item("manufacture")
$ws3.Range("E7").Clear()

# --- workbook-level calculation setting ------------------------------------
$excel.Iteration = $false

# "energy" becomes the active sheet/tab, with E9 (the new row) selected;
# "combustion-co2" (previously active) naturally loses tabSelected.
$ws1.Activate()
$ws1.Range("E9").Select()
